# Auto-generated Excel COM-interop edit script
# Updates the cryptos price list: per-row Price (D) and Volume(1h) (E) refresh,
# plus a ranking reshuffle for rows 43-46 (Coin/Link/Price/Volume all change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.300.60"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "1.903.05"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.37"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4777"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2857"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06679"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.72"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "102.42"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -5.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07707"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "1.911.94"
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.211"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -4.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6737"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "260.16"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -7.97%  "
$ws.Range("D17").Value = "30.350.42"
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007483"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.70"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -4.49%  "
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.281"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.456"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.05"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.91"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.065"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -6.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1009"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.591"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.507"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.214"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04766"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7299"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01918"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.586"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.283"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.47"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -4.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.992"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -6.00%  "
$ws.Range("D47").Value = "1.001.18"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.470"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -8.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1188"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -5.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.862"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -5.82%  "

# Rows 43-46: coin ranking order changed, so Coin/Link/Price/Volume all update
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8595"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.17%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.90"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4244"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -4.87%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.07%  "
